# Error Calculations and Plots
#
# Two fully-missing-data rows are dropped from the dataset ("RM 232" and
# "SC 92"), which shifts every row below them up and shrinks the used
# range from A1:F35 to A1:F33. After the rows shift, a few previously
# missing values get filled in / recomputed:
#   - "SC 5"   (now row 26): column "A" (B26) was blank -> -20.2
#   - "SC 101" (now row 27): column "A" (B27) had -20.4 -> now blank
#   - "SC 232" (now row 33): column "D" (E33) was blank -> -10.7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "RM 232" row (row 26) entirely - remaining rows shift up.
$ws.Rows("26:26").Delete()

# After the shift above, "SC 92" (originally row 28) is now row 27.
# Drop it too - remaining rows shift up again.
$ws.Rows("27:27").Delete()

# Fill in / update values on the rows that shifted into their final
# positions (rows 26-33 after both deletions).
$ws.Range("B26").Value = -20.2
$ws.Range("B27").ClearContents()
$ws.Range("E33").Value = -10.7
